# October month commit: update manager/password credential pair in row 4
# and move the active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "mngr353217"
$ws.Range("B4").Value = "ehadEru"

$ws.Range("B4").Select()
